$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.207.49"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.668.57"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.Value = "'532.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.96%  "
$c = $ws.Range("D6")
$c.Value = "'156.69"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.60%  "
$c = $ws.Range("D7")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$c = $ws.Range("D8")
$c.Value = "'0.591"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "
$c = $ws.Range("D10")
$c.Value = "'0.109"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "3.135.77"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "61.073.88"
$ws.Range("E14").Value = "  +0.88%  "
$c = $ws.Range("D15")
$c.Value = "'22.18"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "2.679.59"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +0.81%  "
$c = $ws.Range("D19")
$c.Value = "'356.34"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("E21").Value = "  +2.63%  "
$c = $ws.Range("D22")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +1.41%  "
$c = $ws.Range("D24")
$c.Value = "'0.434"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D27").Value = "0.0₃0868"
$ws.Range("E27").Value = "  +3.00%  "
$c = $ws.Range("D28")
$c.Value = "'7.43"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("E29").Value = "  -0.02%  "
$c = $ws.Range("D30")
$c.Value = "'6.21"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +6.43%  "
$c = $ws.Range("D31")
$c.Value = "'19.63"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.08%  "
$c = $ws.Range("D32")
$c.Value = "'1.63"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.97%  "
$c = $ws.Range("D33")
$c.Value = "'150.21"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("E34").Value = "  +4.77%  "
$ws.Range("E35").Value = "  +1.46%  "
$c = $ws.Range("D36")
$c.Value = "'0.916"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +8.66%  "
$c = $ws.Range("D37")
$c.Value = "'0.884"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D38")
$c.Value = "'1.51"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D39")
$c.Value = "'309.63"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +5.99%  "
$c = $ws.Range("D40")
$c.Value = "'3.82"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("E42").Value = "  +1.15%  "
$c = $ws.Range("D43")
$c.Value = "'20.61"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D45")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D46")
$c.Value = "'5.04"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("E48").Value = "  +0.47%  "
$c = $ws.Range("D49")
$c.Value = "'19.13"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +8.61%  "
$ws.Range("D50").Value = "2.008.15"
$ws.Range("E50").Value = "  +0.30%  "
$c = $ws.Range("D51")
$c.Value = "'1.85"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.84%  "
